$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content so the sheet ends up with exactly the target cells
$ws.Cells.Clear()

# Row 1
$ws.Range("A1").Value = "Date and Time"
$ws.Range("B1").Value = "2024-03-11 14:55:04.329000 to 2024-03-11 15:49:13.490000"

# Row 2
$ws.Range("A2").Value = "Total time taken for the ride"
$ws.Range("B2").Value = 0.03762708333333333
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

# Row 3
$ws.Range("A3").Value = "Actual Ampere-hours (Ah)"
$ws.Range("B3").Value = 33.50618638888889

# Row 4
$ws.Range("A4").Value = "Actual Watt-hours (Wh)"
$ws.Range("B4").Value = 1709.0979718525

# Row 5
$ws.Range("A5").Value = "Starting SoC (Ah)"
$ws.Range("B5").Value = 37.269

# Row 6
$ws.Range("A6").Value = "Ending SoC (Ah)"
$ws.Range("B6").Value = 3.638

# Row 7
$ws.Range("A7").Value = "Starting SoC (%)"
$ws.Range("B7").Value = 94

# Row 8
$ws.Range("A8").Value = "Ending SoC (%)"
$ws.Range("B8").Value = 9

# Row 9
$ws.Range("A9").Value = "Total distance covered (km)"
$ws.Range("B9").Value = 37.58353708281574

# Row 10
$ws.Range("A10").Value = "Total energy consumption(WH/KM)"
$ws.Range("B10").Value = 45.47464407318778

# Row 11
$ws.Range("A11").Value = "Total SOC consumed(%)"
$ws.Range("B11").Value = 85

# Row 12
$ws.Range("A12").Value = "Mode"
$ws.Range("B12").Value = "Custom mode`n67.79%`nEco mode`n32.21%"

# Row 13
$ws.Range("A13").Value = "Peak Power(kW)"
$ws.Range("B13").Value = 4552.1056

# Row 14
$ws.Range("A14").Value = "Average Power(kW)"
$ws.Range("B14").Value = -1900.75770734291

# Row 15
$ws.Range("A15").Value = "Total Energy Regenerated(kWh)"
$ws.Range("B15").Value = 2.388511011944444

# Row 16
$ws.Range("A16").Value = "Regenerative Effectiveness(%)"
$ws.Range("B16").Value = 0.1395576906892593

# Row 17
$ws.Range("A17").Value = "Highest Cell Voltage(V)"
$ws.Range("B17").Value = 3.349

# Row 18
$ws.Range("A18").Value = "Lowest Cell Voltage(V)"
$ws.Range("B18").Value = 3.047

# Row 19
$ws.Range("A19").Value = "Difference in Cell Voltage(V)"
$ws.Range("B19").Value = 0.302

# Row 20
$ws.Range("A20").Value = "Minimum Temperature(C)"
$ws.Range("B20").Value = 37

# Row 21
$ws.Range("A21").Value = "Maximum Temperature(C)"
$ws.Range("B21").Value = 47

# Row 22
$ws.Range("A22").Value = "Difference in Temperature(C)"
$ws.Range("B22").Value = 10

# Row 23
$ws.Range("A23").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("B23").Value = 70

# Row 24
$ws.Range("A24").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("B24").Value = 67

# Row 25
$ws.Range("A25").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("B25").Value = 65

# Row 26
$ws.Range("A26").Value = "Maximum MCU Temperature(C)"
$ws.Range("B26").Value = 68

# Row 27
$ws.Range("A27").Value = "Maximum Motor Temperature(C)"
$ws.Range("B27").Value = 100

# Row 28
$ws.Range("A28").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("B28").Value = 0

# Row 29
$ws.Range("A29").Value = "highest cell temp(C)"
$ws.Range("B29").Value = 47

# Row 30
$ws.Range("A30").Value = "lowest cell temp(C)"
$ws.Range("B30").Value = 37

# Row 31
$ws.Range("A31").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("B31").Value = 10

# Row 32
$ws.Range("A32").Value = "Battery Voltage(V)"
$ws.Range("B32").Value = 53

# Row 33
$ws.Range("A33").Value = "Total energy charged(kWh)"
$ws.Range("B33").Value = 1.775827878611111

# Row 34
$ws.Range("A34").Value = "Electricity consumption units(kW)"
$ws.Range("B34").Value = 0.0000001517801605650523

# Row 35
$ws.Range("A35").Value = "Cycle Count of battery"
$ws.Range("B35").Value = 51

# Row 36
$ws.Range("A36").Value = "Idling time percentage"
$ws.Range("B36").Value = 0.8162224206096161

# Row 37
$ws.Range("A37").Value = "Time spent in 0-10 km/h"
$ws.Range("B37").Value = 2.79301109552353

# Row 38
$ws.Range("A38").Value = "Time spent in 10-20 km/h"
$ws.Range("B38").Value = 3.75271011350593

# Row 39
$ws.Range("A39").Value = "Time spent in 20-30 km/h"
$ws.Range("B39").Value = 12.41869659482209

# Row 40
$ws.Range("A40").Value = "Time spent in 30-40 km/h"
$ws.Range("B40").Value = 37.85550312460146

# Row 41
$ws.Range("A41").Value = "Time spent in 40-50 km/h"
$ws.Range("B41").Value = 13.14883305700804

# Row 42
$ws.Range("A42").Value = "Time spent in 50-60 km/h"
$ws.Range("B42").Value = 17.48820303532713

# Row 43
$ws.Range("A43").Value = "Time spent in 60-70 km/h"
$ws.Range("B43").Value = 9.217574288993751

# Row 44
$ws.Range("A44").Value = "Time spent in 70-80 km/h"
$ws.Range("B44").Value = 2.426348680015304

# Row 45
$ws.Range("A45").Value = "Time spent in 80-90 km/h"
$ws.Range("B45").Value = 0
